$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.985.17'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').Value = '2.218.33'
$ws.Range('E3').Value = '  -0.20%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '291.98'
$ws.Range('E5').Value = '  -0.14%  '
$ws.Range('D6').Value = '86.96'
$ws.Range('E6').Value = '  -0.76%  '
$ws.Range('E7').Value = '  -0.49%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').Value = '0.468'
$ws.Range('E9').Value = '  -1.14%  '
$ws.Range('D10').Value = '30.44'
$ws.Range('E10').Value = '  +0.28%  '
$ws.Range('D11').Value = '50.31'
$ws.Range('E11').Value = '  +5.93%  '
$ws.Range('D12').Value = '0.0779'
$ws.Range('E12').Value = '  -0.98%  '
$ws.Range('E13').Value = '  +3.06%  '
$ws.Range('D14').Value = '6.44'
$ws.Range('E14').Value = '  +0.20%  '
$ws.Range('D15').Value = '2.558.58'
$ws.Range('E15').Value = '  -0.17%  '
$ws.Range('D16').Value = '13.78'
$ws.Range('E16').Value = '  -2.34%  '
$ws.Range('D17').Value = '2.249.76'
$ws.Range('E17').Value = '  +1.25%  '
$ws.Range('D18').Value = '0.732'
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('D19').Value = '39.879.57'
$ws.Range('E19').Value = '  +0.11%  '
$ws.Range('D20').Value = '0.0₃0886'
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('D21').Value = '11.14'
$ws.Range('E21').Value = '  -3.24%  '
$ws.Range('D22').Value = '5.75'
$ws.Range('E22').Value = '  -1.63%  '
$ws.Range('D23').Value = '65.65'
$ws.Range('E23').Value = '  -0.30%  '
$ws.Range('D24').Value = '237.88'
$ws.Range('E24').Value = '  +0.84%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('D26').Value = '2.45'
$ws.Range('E26').Value = '  -0.96%  '
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('D28').Value = '23.10'
$ws.Range('E28').Value = '  +1.11%  '
$ws.Range('D29').Value = '9.24'
$ws.Range('E29').Value = '  -0.46%  '
$ws.Range('E30').Value = '  -7.72%  '
$ws.Range('D31').Value = '156.79'
$ws.Range('E31').Value = '  +2.80%  '
$ws.Range('D32').Value = '31.88'
$ws.Range('E32').Value = '  -2.97%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').Value = '4.96'
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('D35').Value = '2.99'
$ws.Range('E35').Value = '  +6.03%  '
$ws.Range('D36').Value = '0.0715'
$ws.Range('E36').Value = '  -0.86%  '
$ws.Range('E37').Value = '  -1.67%  '
$ws.Range('E38').Value = '  -0.14%  '
$ws.Range('D39').Value = '0.0992'
$ws.Range('E39').Value = '  -0.14%  '
$ws.Range('E40').Value = '  +0.89%  '
$ws.Range('D41').Value = '15.26'
$ws.Range('E41').Value = '  -4.55%  '
$ws.Range('D42').Value = '2.093.69'
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('D43').Value = '3.71'
$ws.Range('E43').Value = '  -2.73%  '
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('D45').Value = '18.09'
$ws.Range('E45').Value = '  +1.69%  '
$ws.Range('D46').Value = '9.79'
$ws.Range('E46').Value = '  -2.17%  '
$ws.Range('E47').Value = '  -8.11%  '
$ws.Range('D48').Value = '2.71'
$ws.Range('E48').Value = '  +2.57%  '
$ws.Range('D49').Value = '2.433.04'
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('D50').Value = '1.46'
$ws.Range('E50').Value = '  +0.23%  '
$ws.Range('D51').Value = '1.11'
$ws.Range('E51').Value = '  +2.37%  '
